# #5: fund, bonds, otherbonds, antique done
#
# The "基金受益憑證" (fund) sheet (sheet index 4) previously had a bogus
# header row 1 that was just a duplicate of row 2's data. This fixes row 1
# to be a real header row, and appends the standard metadata columns
# (property_category .. index) that every other sheet in this workbook
# already carries, to every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# ---- Row 1: real header labels (was a stray duplicate of row 2) ----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# ---- Rows 2-11: append the standard metadata columns I:O ----
# property_category = fund, category = normal, date = 2013-12-12,
# legislator_name = 吳秉叡, legislator_id = 1324, source_file = tmpea101,
# index = same value as column A on that row.
$indices = @{ 2 = 90; 3 = 92; 4 = 93; 5 = 94; 6 = 95; 7 = 96; 8 = 97; 9 = 98; 10 = 99; 11 = 100 }

foreach ($row in 2..11) {
    $ws.Range("I$row").Value = "fund"
    $ws.Range("J$row").Value = "normal"
    $ws.Range("K$row").Value = "2013-12-12"
    $ws.Range("L$row").Value = "吳秉叡"
    $ws.Range("M$row").Value = 1324
    $ws.Range("N$row").Value = "tmpea101"
    $ws.Range("O$row").Value = $indices[$row]
}
